# Updated symbol list (crypto prices / 1h volume changes) as scraped by
# GitHub Actions. Values are written with a leading apostrophe (quote
# prefix) via .Formula so Excel stores them as literal text (matching the
# original inlineStr cells) instead of auto-converting to numbers/percentages.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Formula = "'274.38"
$ws.Range("E2").Formula = "'-1.24%"
$ws.Range("D3").Formula = "'27.07"
$ws.Range("E3").Formula = "'-1.29%"
$ws.Range("D4").Formula = "'4.790"
$ws.Range("E4").Formula = "'-0.44%"
$ws.Range("D5").Formula = "'0.06283"
$ws.Range("E5").Formula = "'-0.70%"
$ws.Range("D6").Formula = "'6.914"
$ws.Range("E6").Formula = "'-0.28%"
$ws.Range("D7").Formula = "'1.306"
$ws.Range("E7").Formula = "'38.22%"
$ws.Range("D8").Formula = "'0.8702"
$ws.Range("E8").Formula = "'-1.59%"
$ws.Range("D9").Formula = "'0.1554"
$ws.Range("E9").Formula = "'5.16%"
$ws.Range("D10").Formula = "'0.05003"
$ws.Range("E10").Formula = "'-4.90%"
$ws.Range("D11").Formula = "'0.07481"
$ws.Range("E11").Formula = "'1.80%"
$ws.Range("D12").Formula = "'0.02883"
$ws.Range("E12").Formula = "'-8.52%"
$ws.Range("D13").Formula = "'0.09055"
$ws.Range("D14").Formula = "'0.001568"
$ws.Range("E14").Formula = "'1.20%"
$ws.Range("D15").Formula = "'0.0006356"
$ws.Range("E15").Formula = "'1.24%"
$ws.Range("D16").Formula = "'0.005874"
$ws.Range("E16").Formula = "'1.31%"
$ws.Range("D17").Formula = "'3.453"
$ws.Range("E17").Formula = "'-0.36%"
$ws.Range("D18").Formula = "'3.310"
$ws.Range("E18").Formula = "'-1.44%"
$ws.Range("E19").Formula = "'-0.44%"
$ws.Range("D20").Formula = "'0.3148"
$ws.Range("E20").Formula = "'1.72%"
$ws.Range("D21").Formula = "'0.1318"
$ws.Range("E21").Formula = "'-1.48%"
$ws.Range("D22").Formula = "'3.925"
$ws.Range("E22").Formula = "'0.98%"
$ws.Range("D23").Formula = "'0.04395"
$ws.Range("E23").Formula = "'1.69%"
$ws.Range("E24").Formula = "'-1.06%"
$ws.Range("D25").Formula = "'0.003811"
$ws.Range("E25").Formula = "'6.27%"
$ws.Range("E26").Formula = "'0.14%"
$ws.Range("D27").Formula = "'0.0001617"
$ws.Range("E27").Formula = "'-4.51%"
$ws.Range("D40").Formula = "'0.04073"
$ws.Range("E40").Formula = "'0.59%"
$ws.Range("D41").Formula = "'0.007062"
$ws.Range("E41").Formula = "'5.82%"
$ws.Range("E42").Formula = "'0.31%"
$ws.Range("D43").Formula = "'0.002032"
$ws.Range("E43").Formula = "'-12.58%"
$ws.Range("E44").Formula = "'-9.83%"
$ws.Range("D45").Formula = "'0.00005176"
$ws.Range("E45").Formula = "'-1.22%"
$ws.Range("D46").Formula = "'0.02302"
$ws.Range("E46").Formula = "'2.11%"
$ws.Range("D47").Formula = "'1.490"
$ws.Range("E47").Formula = "'-37.34%"
